# Refresh the crypto price table with the latest scraped values.
# Price-looking text in column D (e.g. "245.24") must remain plain text,
# exactly like the source data (Excel would otherwise coerce it to a number),
# so those assignments are written with a leading apostrophe, matching how
# Excel stores manually quote-prefixed text values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '36.971.53'
$ws.Range('E2').Value = '  -0.56%  '
# Row 3
$ws.Range('D3').Value = '2.038.84'
$ws.Range('E3').Value = '  -0.68%  '
# Row 4
$ws.Range('E4').Value = '  -0.02%  '
# Row 5
$ws.Range('D5').Value = '''245.24'
$ws.Range('E5').Value = '  -1.86%  '
# Row 6
$ws.Range('E6').Value = '  -1.75%  '
# Row 7
$ws.Range('D7').Value = '''58.39'
$ws.Range('E7').Value = '  -1.61%  '
# Row 8
$ws.Range('E8').Value = '  -0.02%  '
# Row 9
$ws.Range('E9').Value = '  -2.15%  '
# Row 10
$ws.Range('E10').Value = '  -2.25%  '
# Row 11
$ws.Range('E11').Value = '  +2.28%  '
# Row 12
$ws.Range('D12').Value = '''15.30'
$ws.Range('E12').Value = '  -5.53%  '
# Row 13
$ws.Range('D13').Value = '''0.882'
$ws.Range('E13').Value = '  +7.99%  '
# Row 14
$ws.Range('D14').Value = '2.335.18'
$ws.Range('E14').Value = '  -0.71%  '
# Row 15
$ws.Range('D15').Value = '''5.63'
$ws.Range('E15').Value = '  +0.59%  '
# Row 16
$ws.Range('D16').Value = '2.025.48'
$ws.Range('E16').Value = '  -1.33%  '
# Row 17
$ws.Range('D17').Value = '''18.21'
$ws.Range('E17').Value = '  +2.58%  '
# Row 18
$ws.Range('D18').Value = '36.934.00'
$ws.Range('E18').Value = '  -0.63%  '
# Row 19
$ws.Range('D19').Value = '''73.44'
$ws.Range('E19').Value = '  -1.94%  '
# Row 20
$ws.Range('E20').Value = '  -2.34%  '
# Row 21
$ws.Range('E21').Value = '  -0.63%  '
# Row 22
$ws.Range('D22').Value = '''235.04'
$ws.Range('E22').Value = '  -1.07%  '
# Row 23
$ws.Range('E23').Value = '  +0.03%  '
# Row 24
$ws.Range('E24').Value = '  +1.75%  '
# Row 25
$ws.Range('D25').Value = '''9.54'
$ws.Range('E25').Value = '  +1.94%  '
# Row 26
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '''167.83'
$ws.Range('E26').Value = '  -0.45%  '
# Row 27
$ws.Range('B27').Value = 'PancakeSwap'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D27').Value = '''2.14'
$ws.Range('E27').Value = '  -2.57%  '
# Row 28
$ws.Range('D28').Value = '''19.87'
$ws.Range('E28').Value = '  -0.43%  '
# Row 29
$ws.Range('D29').Value = '''5.54'
$ws.Range('E29').Value = '  +16.28%  '
# Row 31
$ws.Range('E31').Value = '  -3.56%  '
# Row 32
$ws.Range('D32').Value = '''4.76'
$ws.Range('E32').Value = '  +5.48%  '
# Row 33
$ws.Range('E33').Value = '  -0.97%  '
# Row 34
$ws.Range('E34').Value = '  +0.06%  '
# Row 35
$ws.Range('D35').Value = '''0.0856'
$ws.Range('E35').Value = '  -4.49%  '
# Row 36
$ws.Range('E36').Value = '  +5.67%  '
# Row 37
$ws.Range('E37').Value = '  -0.56%  '
# Row 38
$ws.Range('E38').Value = '  -4.87%  '
# Row 39
$ws.Range('D39').Value = '''5.23'
$ws.Range('E39').Value = '  -1.43%  '
# Row 40
$ws.Range('E40').Value = '  -1.71%  '
# Row 41
$ws.Range('E41').Value = '  -0.26%  '
# Row 42
$ws.Range('E42').Value = '  +0.67%  '
# Row 43
$ws.Range('D43').Value = '''0.0946'
$ws.Range('E43').Value = '  -14.74%  '
# Row 44
$ws.Range('D44').Value = '''96.91'
$ws.Range('E44').Value = '  +0.65%  '
# Row 45
$ws.Range('D45').Value = '''16.83'
$ws.Range('E45').Value = '  -4.98%  '
# Row 46
$ws.Range('D46').Value = '1.290.93'
$ws.Range('E46').Value = '  +0.37%  '
# Row 47
$ws.Range('D47').Value = '''2.36'
$ws.Range('E47').Value = '  -4.43%  '
# Row 48
$ws.Range('E48').Value = '  -0.64%  '
# Row 49
$ws.Range('B49').Value = 'FTXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D49').Value = '''3.63'
$ws.Range('E49').Value = '  +6.66%  '
# Row 50
$ws.Range('B50').Value = 'FraxShare'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D50').Value = '''6.66'
$ws.Range('E50').Value = '  -2.04%  '
# Row 51
$ws.Range('D51').Value = '2.220.70'
$ws.Range('E51').Value = '  -1.05%  '
